# Auto-generated COM-interop script applying the 2026-01-28 06:54:14 scrape refresh
# to the three schedule sheets (LP1912, LP1912-215, 6203-6173).
$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912  (121 cell changes) ----
$ws = $wb.Worksheets.Item("LP1912")
$ws.Cells.Item(2,1).Value = 'Última actualización: 06:54:14'
$ws.Cells.Item(3,1).Value = 'Total filas: 65'
$ws.Cells.Item(30,3).Value = '86_EST CHICA-ESC AGRARIA'
$ws.Cells.Item(31,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(38,1).Value = '06:54:14'
$ws.Cells.Item(38,4).Value = 6
$ws.Cells.Item(41,1).Value = '06:54:14'
$ws.Cells.Item(41,4).Value = 11
$ws.Cells.Item(42,1).Value = '06:54:14'
$ws.Cells.Item(42,4).Value = 13
$ws.Cells.Item(43,1).Value = '06:54:14'
$ws.Cells.Item(43,4).Value = 17
$ws.Cells.Item(44,1).Value = '06:54:14'
$ws.Cells.Item(44,2).Value = '07:11'
$ws.Cells.Item(44,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(44,4).Value = 17
$ws.Cells.Item(45,1).Value = '06:26:08'
$ws.Cells.Item(45,2).Value = '07:12'
$ws.Cells.Item(45,3).Value = '215A_EL PATO'
$ws.Cells.Item(45,4).Value = 46
$ws.Cells.Item(46,1).Value = '05:55:25'
$ws.Cells.Item(46,2).Value = '07:15'
$ws.Cells.Item(46,4).Value = 80
$ws.Cells.Item(47,1).Value = '06:54:14'
$ws.Cells.Item(47,2).Value = '07:16'
$ws.Cells.Item(47,3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(47,4).Value = 22
$ws.Cells.Item(48,1).Value = '06:54:14'
$ws.Cells.Item(48,2).Value = '07:17'
$ws.Cells.Item(48,3).Value = '16_SANTA ANA'
$ws.Cells.Item(48,4).Value = 23
$ws.Cells.Item(49,1).Value = '06:54:14'
$ws.Cells.Item(49,2).Value = '07:21'
$ws.Cells.Item(49,3).Value = '26_HERNANDEZ'
$ws.Cells.Item(49,4).Value = 27
$ws.Cells.Item(50,1).Value = '06:54:14'
$ws.Cells.Item(50,2).Value = '07:23'
$ws.Cells.Item(50,3).Value = '10_OLMOS'
$ws.Cells.Item(50,4).Value = 29
$ws.Cells.Item(51,1).Value = '05:55:25'
$ws.Cells.Item(51,2).Value = '07:31'
$ws.Cells.Item(51,3).Value = '16_SANTA ANA'
$ws.Cells.Item(51,4).Value = 96
$ws.Cells.Item(52,1).Value = '05:55:25'
$ws.Cells.Item(52,2).Value = '07:31'
$ws.Cells.Item(52,4).Value = 96
$ws.Cells.Item(53,1).Value = '06:54:14'
$ws.Cells.Item(53,2).Value = '07:32'
$ws.Cells.Item(53,3).Value = '16_SANTA ANA'
$ws.Cells.Item(53,4).Value = 38
$ws.Cells.Item(54,1).Value = '06:54:14'
$ws.Cells.Item(54,2).Value = '07:32'
$ws.Cells.Item(54,3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(54,4).Value = 38
$ws.Cells.Item(55,1).Value = '06:54:14'
$ws.Cells.Item(55,2).Value = '07:32'
$ws.Cells.Item(55,3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws.Cells.Item(55,4).Value = 38
$ws.Cells.Item(56,2).Value = '07:36'
$ws.Cells.Item(56,3).Value = '27_EL RETIRO'
$ws.Cells.Item(56,4).Value = 101
$ws.Cells.Item(57,1).Value = '06:54:14'
$ws.Cells.Item(57,2).Value = '07:37'
$ws.Cells.Item(57,3).Value = '27_EL RETIRO'
$ws.Cells.Item(57,4).Value = 43
$ws.Cells.Item(58,1).Value = '06:54:14'
$ws.Cells.Item(58,2).Value = '07:39'
$ws.Cells.Item(58,3).Value = '10_OLMOS'
$ws.Cells.Item(58,4).Value = 45
$ws.Cells.Item(59,1).Value = '05:55:25'
$ws.Cells.Item(59,2).Value = '07:47'
$ws.Cells.Item(59,3).Value = '14_ABASTO'
$ws.Cells.Item(59,4).Value = 112
$ws.Cells.Item(60,1).Value = '06:54:14'
$ws.Cells.Item(60,2).Value = '07:48'
$ws.Cells.Item(60,3).Value = '14_ABASTO'
$ws.Cells.Item(60,4).Value = 54
$ws.Cells.Item(61,1).Value = '06:54:14'
$ws.Cells.Item(61,2).Value = '07:51'
$ws.Cells.Item(61,3).Value = '215D_EL PATO'
$ws.Cells.Item(61,4).Value = 57
$ws.Cells.Item(62,2).Value = '07:52'
$ws.Cells.Item(62,3).Value = '215D_EL PATO'
$ws.Cells.Item(62,4).Value = 86
$ws.Cells.Item(63,2).Value = '08:01'
$ws.Cells.Item(63,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(63,4).Value = 95
$ws.Cells.Item(64,1).Value = '06:54:14'
$ws.Cells.Item(64,2).Value = '08:06'
$ws.Cells.Item(64,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(64,4).Value = 72
$ws.Cells.Item(65,1).Value = '06:54:14'
$ws.Cells.Item(65,2).Value = '08:12'
$ws.Cells.Item(65,3).Value = '15_ABASTO'
$ws.Cells.Item(65,4).Value = 78
$ws.Cells.Item(65,5).Value = 'LP1912'
$ws.Cells.Item(66,1).Value = '06:54:14'
$ws.Cells.Item(66,2).Value = '08:21'
$ws.Cells.Item(66,3).Value = '26_HERNANDEZ'
$ws.Cells.Item(66,4).Value = 87
$ws.Cells.Item(66,5).Value = 'LP1912'
$ws.Cells.Item(67,1).Value = '06:54:14'
$ws.Cells.Item(67,2).Value = '08:23'
$ws.Cells.Item(67,3).Value = '215B_EL PATO'
$ws.Cells.Item(67,4).Value = 89
$ws.Cells.Item(67,5).Value = 'LP1912'
$ws.Cells.Item(68,1).Value = '06:54:14'
$ws.Cells.Item(68,2).Value = '08:23'
$ws.Cells.Item(68,3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(68,4).Value = 89
$ws.Cells.Item(68,5).Value = 'LP1912'
$ws.Cells.Item(69,1).Value = '06:54:14'
$ws.Cells.Item(69,2).Value = '08:27'
$ws.Cells.Item(69,3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws.Cells.Item(69,4).Value = 93
$ws.Cells.Item(69,5).Value = 'LP1912'
$ws.Cells.Item(70,1).Value = '06:54:14'
$ws.Cells.Item(70,2).Value = '08:42'
$ws.Cells.Item(70,3).Value = '81_EL PELIGRO'
$ws.Cells.Item(70,4).Value = 108
$ws.Cells.Item(70,5).Value = 'LP1912'

# ---- Sheet: LP1912-215  (7 cell changes) ----
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Cells.Item(2,1).Value = 'Última actualización: 06:54:14'
$ws.Cells.Item(14,1).Value = '06:54:14'
$ws.Cells.Item(14,4).Value = 17
$ws.Cells.Item(16,1).Value = '06:54:14'
$ws.Cells.Item(16,4).Value = 57
$ws.Cells.Item(18,1).Value = '06:54:14'
$ws.Cells.Item(18,4).Value = 89

# ---- Sheet: 6203-6173  (13 cell changes) ----
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Cells.Item(2,1).Value = 'Última actualización: 06:54:14'
$ws.Cells.Item(3,1).Value = 'Total filas: 10'
$ws.Cells.Item(12,1).Value = '06:54:14'
$ws.Cells.Item(12,4).Value = 6
$ws.Cells.Item(13,1).Value = '06:54:14'
$ws.Cells.Item(13,4).Value = 41
$ws.Cells.Item(14,1).Value = '06:54:14'
$ws.Cells.Item(14,4).Value = 73
$ws.Cells.Item(15,1).Value = '06:54:14'
$ws.Cells.Item(15,2).Value = '08:31'
$ws.Cells.Item(15,3).Value = '215A_LA PLATA'
$ws.Cells.Item(15,4).Value = 97
$ws.Cells.Item(15,5).Value = 'L6173'

Write-Output "Schedules updated."
